# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list ... with GitHub Actions" (Tue May  2 20:53:40 UTC 2023).
# Columns B-E hold plain text (coin name / link / price / 1h volume), so we
# force text format on D/E before writing to stop Excel from re-parsing
# strings like "0.9800" or "17.00" as numbers and dropping trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.848.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.881.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.45%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.04"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4673"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3947"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.68%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07934"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9800"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.34"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.80%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.895.81"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.750"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.019"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06966"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.68"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.007"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001010"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.91%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.00"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.32%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.871.50"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.369"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.12"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.125"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.114.31"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.60"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.44"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.772"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.007"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.09"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09407"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.87%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9426"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.323"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.357"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.354"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05921"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.59%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.926"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.97%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5725"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.10%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1797"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.97%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.01"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07290"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.86"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5350"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.36%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.86%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.118"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.852"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "114.45"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.373"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.44%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.46%  "
